$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 18547
$ws.Range("B4").Value = 12020
$ws.Range("B5").Value = 30567
$ws.Range("B6").Value = 0.14338375
